# Apply targeted cell updates per the commit diff.
# Stricter separation of department specific pre mid and post mid courses:
# room/facility/capacity reassignments on Classroom_Allocation, and
# updated "Allocated Rooms" summaries on Basket_Course_Allocations.

$wb = $excel.ActiveWorkbook

# --- Sheet: Classroom_Allocation ---
$ws6 = $wb.Worksheets.Item("Classroom_Allocation")

$ws6.Range("I5").Value = "Projector"
$ws6.Range("M5").Value = "C002"
$ws6.Range("G6").Value = "classroom"
$ws6.Range("H6").NumberFormat = "@"
$ws6.Range("H6").Value = "96"
$ws6.Range("H6").Style = "Normal"
$ws6.Range("M6").Value = "C102"
$ws6.Range("M7").Value = "C104"
$ws6.Range("I8").Value = "TV"
$ws6.Range("M8").Value = "C203"
$ws6.Range("I9").Value = "TV"
$ws6.Range("M9").Value = "C205"
$ws6.Range("G17").Value = "Auditorium"
$ws6.Range("H17").NumberFormat = "@"
$ws6.Range("H17").Value = "240"
$ws6.Range("H17").Style = "Normal"
$ws6.Range("I17").Value = "Audio/Video System"
$ws6.Range("M17").Value = "C004"
$ws6.Range("G18").Value = "large classroom"
$ws6.Range("H18").NumberFormat = "@"
$ws6.Range("H18").Value = "120"
$ws6.Range("H18").Style = "Normal"
$ws6.Range("I18").Value = ""
$ws6.Range("M18").Value = "C001"
$ws6.Range("G19").Value = "large classroom"
$ws6.Range("H19").NumberFormat = "@"
$ws6.Range("H19").Value = "120"
$ws6.Range("H19").Style = "Normal"
$ws6.Range("M19").Value = "C002"
$ws6.Range("M20").Value = "C101"
$ws6.Range("I21").Value = "Projector"
$ws6.Range("M21").Value = "C002"
$ws6.Range("G22").Value = "classroom"
$ws6.Range("H22").NumberFormat = "@"
$ws6.Range("H22").Value = "96"
$ws6.Range("H22").Style = "Normal"
$ws6.Range("M22").Value = "C102"
$ws6.Range("M23").Value = "C104"
$ws6.Range("I24").Value = "TV"
$ws6.Range("M24").Value = "C203"
$ws6.Range("I25").Value = "TV"
$ws6.Range("M25").Value = "C205"
$ws6.Range("G37").Value = "large classroom"
$ws6.Range("H37").NumberFormat = "@"
$ws6.Range("H37").Value = "120"
$ws6.Range("H37").Style = "Normal"
$ws6.Range("M37").Value = "C002"
$ws6.Range("M38").Value = "C102"
$ws6.Range("M39").Value = "C104"
$ws6.Range("M41").Value = "C205"
$ws6.Range("G42").Value = "Auditorium"
$ws6.Range("H42").NumberFormat = "@"
$ws6.Range("H42").Value = "240"
$ws6.Range("H42").Style = "Normal"
$ws6.Range("I42").Value = "Audio/Video System"
$ws6.Range("M42").Value = "C004"
$ws6.Range("G43").Value = "large classroom"
$ws6.Range("H43").NumberFormat = "@"
$ws6.Range("H43").Value = "120"
$ws6.Range("H43").Style = "Normal"
$ws6.Range("I43").Value = ""
$ws6.Range("M43").Value = "C001"
$ws6.Range("G44").Value = "large classroom"
$ws6.Range("H44").NumberFormat = "@"
$ws6.Range("H44").Value = "120"
$ws6.Range("H44").Style = "Normal"
$ws6.Range("M44").Value = "C002"
$ws6.Range("M45").Value = "C101"
$ws6.Range("G46").Value = "Auditorium"
$ws6.Range("H46").NumberFormat = "@"
$ws6.Range("H46").Value = "240"
$ws6.Range("H46").Style = "Normal"
$ws6.Range("I46").Value = "Audio/Video System"
$ws6.Range("M46").Value = "C004"
$ws6.Range("G47").Value = "large classroom"
$ws6.Range("H47").NumberFormat = "@"
$ws6.Range("H47").Value = "120"
$ws6.Range("H47").Style = "Normal"
$ws6.Range("I47").Value = ""
$ws6.Range("M47").Value = "C001"
$ws6.Range("G48").Value = "large classroom"
$ws6.Range("H48").NumberFormat = "@"
$ws6.Range("H48").Value = "120"
$ws6.Range("H48").Style = "Normal"
$ws6.Range("M48").Value = "C002"
$ws6.Range("I49").Value = "Projector"
$ws6.Range("M49").Value = "C101"
$ws6.Range("G54").Value = "large classroom"
$ws6.Range("H54").NumberFormat = "@"
$ws6.Range("H54").Value = "120"
$ws6.Range("H54").Style = "Normal"
$ws6.Range("I54").Value = ""
$ws6.Range("M54").Value = "C001"
$ws6.Range("I55").Value = "Projector"
$ws6.Range("M55").Value = "C101"
$ws6.Range("I56").Value = "Projector"
$ws6.Range("M56").Value = "C202"
$ws6.Range("M57").Value = "C204"
$ws6.Range("G65").Value = "classroom"
$ws6.Range("H65").NumberFormat = "@"
$ws6.Range("H65").Value = "96"
$ws6.Range("H65").Style = "Normal"
$ws6.Range("I65").Value = "Projector"
$ws6.Range("M65").Value = "C102"
$ws6.Range("G66").Value = "classroom"
$ws6.Range("H66").NumberFormat = "@"
$ws6.Range("H66").Value = "96"
$ws6.Range("H66").Style = "Normal"
$ws6.Range("I66").Value = "Projector"
$ws6.Range("M66").Value = "C104"
$ws6.Range("G70").Value = "large classroom"
$ws6.Range("H70").NumberFormat = "@"
$ws6.Range("H70").Value = "120"
$ws6.Range("H70").Style = "Normal"
$ws6.Range("I70").Value = ""
$ws6.Range("M70").Value = "C001"
$ws6.Range("I71").Value = "Projector"
$ws6.Range("M71").Value = "C101"
$ws6.Range("I72").Value = "Projector"
$ws6.Range("M72").Value = "C202"
$ws6.Range("M73").Value = "C204"
$ws6.Range("G87").Value = "classroom"
$ws6.Range("H87").NumberFormat = "@"
$ws6.Range("H87").Value = "96"
$ws6.Range("H87").Style = "Normal"
$ws6.Range("M87").Value = "C101"
$ws6.Range("M88").Value = "C202"
$ws6.Range("M89").Value = "C204"
$ws6.Range("G90").Value = "classroom"
$ws6.Range("H90").NumberFormat = "@"
$ws6.Range("H90").Value = "96"
$ws6.Range("H90").Style = "Normal"
$ws6.Range("I90").Value = "Projector"
$ws6.Range("M90").Value = "C102"
$ws6.Range("G91").Value = "classroom"
$ws6.Range("H91").NumberFormat = "@"
$ws6.Range("H91").Value = "96"
$ws6.Range("H91").Style = "Normal"
$ws6.Range("I91").Value = "Projector"
$ws6.Range("M91").Value = "C104"
$ws6.Range("G94").Value = "classroom"
$ws6.Range("H94").NumberFormat = "@"
$ws6.Range("H94").Value = "96"
$ws6.Range("H94").Style = "Normal"
$ws6.Range("I94").Value = "Projector"
$ws6.Range("M94").Value = "C102"
$ws6.Range("G95").Value = "classroom"
$ws6.Range("H95").NumberFormat = "@"
$ws6.Range("H95").Value = "96"
$ws6.Range("H95").Style = "Normal"
$ws6.Range("I95").Value = "Projector"
$ws6.Range("M95").Value = "C104"
$ws6.Range("G96").Value = "classroom"
$ws6.Range("H96").NumberFormat = "@"
$ws6.Range("H96").Value = "96"
$ws6.Range("H96").Style = "Normal"
$ws6.Range("M96").Value = "C202"
$ws6.Range("I97").Value = "TV"
$ws6.Range("M97").Value = "C203"

# --- Sheet: Basket_Course_Allocations ---
$ws7 = $wb.Worksheets.Item("Basket_Course_Allocations")

$ws7.Range("C9").Value = "C002, C004"
$ws7.Range("C10").Value = "C001, C102"
$ws7.Range("C11").Value = "C101, C104"
$ws7.Range("C12").Value = "C202, C203"
$ws7.Range("C13").Value = "C204, C205"
$ws7.Range("C14").Value = "C004, C102"
$ws7.Range("C15").Value = "C001, C104"
$ws7.Range("C16").Value = "C002, C202"
$ws7.Range("C17").Value = "C101, C203"
